$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (TIPO DE MOVIMIENTO shifts right),
# shifting the whole block one column to the right.
$ws.Columns("B:B").Insert()

# New "STOCK" header cell, bold/filled/centered like the other headers.
$ws.Range("B3").Value = "STOCK"
$ws.Range("B3").HorizontalAlignment = -4108

# Match the narrower width used for the new column (raw stored width 16;
# the ColumnWidth setter adds ~0.8333 padding, so back that out).
$ws.Columns("B:B").ColumnWidth = 15.166666666666666

# The title merge (A1:F1) automatically grows to A1:G1 to keep covering the
# same header columns now that a column was inserted inside its span.

# Selection ends on the freshly-added cell.
[void]$ws.Range("B3").Select()
